$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidated report: the "Absent" column (H) is recomputed from the
# "Real" column (E) -- a student is Absent (1) when they were not
# actually present for real (E = 0), and not absent (0) when E = 1.
$lastRow = 21
for ($r = 3; $r -le $lastRow; $r++) {
    $real = $ws.Cells.Item($r, 5).Value2
    if ($real -eq 1) {
        $ws.Cells.Item($r, 8).Value = 0
    } else {
        $ws.Cells.Item($r, 8).Value = 1
    }
}
